$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Application name" column (F) — billing/app-identification detail.
$ws.Range("F1").Value = "Application name"

# Match the bold / centered / wrapped header style already used by B1:D1.
$ws.Range("D1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# TestID 1's row gets the app name value.
$ws.Range("F2").Value = "demowebshop"

# The header row grows to the same row height as the wrapped data rows.
$ws.Rows.Item(1).RowHeight = 28.8

# Size the new column.
$ws.Columns.Item(6).ColumnWidth = 15.17

# Restore the workbook's active selection.
$ws.Range("F7").Select()
